$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metadata")

# Publisher (row 9) - translate from German to English
$ws.Range("B9").Value = "Independent Trusted Third Party of the University Medicine Greifswald"

# Contact (row 10) - translate from German to English
$ws.Range("B10").Value = "Independent Trusted Third Party of the University Medicine Greifswald (https://www.ths-greifswald.de/)"

# Description (row 12) - was empty, now filled in
$ws.Range("B12").Value = "Extensible description of the permissible (data processing-) action resulting from consented module or policy. "
